# Updates cryptos list figures (price + 1h volume change) to match the
# latest scrape, and fixes the swapped Quant/NEARProtocol rows (47 & 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text write - safe for values that Excel would never mis-parse as a
# number/date (percentages with surrounding spaces, coin names, URLs, and
# "price" strings that already contain two dots such as "27.365.63").
function Set-TextValue {
    param($sheet, $cellRef, $val)
    $sheet.Range($cellRef).Value2 = $val
}

# Text write for "Price" column entries that otherwise look like a genuine
# number (e.g. "21.00", "0.4467") - force the cell to Text format first so
# Excel's auto-detection doesn't turn the string into a numeric value, then
# restore the Normal style so we don't leave a stray "@" number format
# behind on the cell.
function Set-TextNumberLike {
    param($sheet, $cellRef, $val)
    $sheet.Range($cellRef).NumberFormat = "@"
    $sheet.Range($cellRef).Value2 = $val
    $sheet.Range($cellRef).Style = "Normal"
}

Set-TextNumberLike $ws "D2" "27.365.63"
Set-TextValue $ws "E2" "  +1.27%  "

Set-TextNumberLike $ws "D3" "1.824.29"
Set-TextValue $ws "E3" "  -0.12%  "

Set-TextValue $ws "E4" "  +0.02%  "

Set-TextValue $ws "E5" "  +0.58%  "

Set-TextValue $ws "E6" "  +0.04%  "

Set-TextNumberLike $ws "D7" "0.4467"
Set-TextValue $ws "E7" "  +2.73%  "

Set-TextValue $ws "E8" "  +2.00%  "

Set-TextNumberLike $ws "D9" "0.07497"
Set-TextValue $ws "E9" "  +3.12%  "

Set-TextNumberLike $ws "D10" "0.8859"
Set-TextValue $ws "E10" "  +4.84%  "

Set-TextNumberLike $ws "D11" "21.00"
Set-TextValue $ws "E11" "  +1.58%  "

Set-TextNumberLike $ws "D12" "1.826.36"
Set-TextValue $ws "E12" "  -0.23%  "

Set-TextNumberLike $ws "D13" "6.752"
Set-TextValue $ws "E13" "  +1.34%  "

Set-TextNumberLike $ws "D14" "93.95"
Set-TextValue $ws "E14" "  +4.88%  "

Set-TextNumberLike $ws "D15" "5.402"
Set-TextValue $ws "E15" "  +1.98%  "

Set-TextNumberLike $ws "D16" "0.07103"
Set-TextValue $ws "E16" "  +0.59%  "

Set-TextValue $ws "E17" "  -0.06%  "

Set-TextNumberLike $ws "D20" "15.22"
Set-TextValue $ws "E20" "  +1.95%  "

Set-TextNumberLike $ws "D21" "27.375.83"
Set-TextValue $ws "E21" "  +0.92%  "

Set-TextNumberLike $ws "D22" "5.256"
Set-TextValue $ws "E22" "  +2.09%  "

Set-TextValue $ws "E23" "  +0.17%  "

Set-TextNumberLike $ws "D24" "1.961"
Set-TextValue $ws "E24" "  -1.61%  "

Set-TextNumberLike $ws "D25" "2.380"
Set-TextValue $ws "E25" "  +7.42%  "

Set-TextNumberLike $ws "D26" "151.51"
Set-TextValue $ws "E26" "  -0.03%  "

Set-TextNumberLike $ws "D27" "18.57"
Set-TextValue $ws "E27" "  +1.43%  "

Set-TextNumberLike $ws "D28" "5.363"
Set-TextValue $ws "E28" "  +2.47%  "

Set-TextNumberLike $ws "D29" "118.02"
Set-TextValue $ws "E29" "  +0.87%  "

Set-TextValue $ws "E30" "  +0.60%  "

Set-TextNumberLike $ws "D31" "0.7841"
Set-TextValue $ws "E31" "  +5.63%  "

Set-TextNumberLike $ws "D32" "1.195"
Set-TextValue $ws "E32" "  +1.27%  "

Set-TextNumberLike $ws "D33" "4.510"
Set-TextValue $ws "E33" "  +1.61%  "

Set-TextNumberLike $ws "D34" "2.931"
Set-TextValue $ws "E34" "  +0.86%  "

Set-TextValue $ws "E35" "  +0.02%  "

Set-TextNumberLike $ws "D36" "1.112"
Set-TextValue $ws "E36" "  +1.35%  "

Set-TextValue $ws "E37" "  +2.30%  "

Set-TextValue $ws "E38" "  +1.68%  "

Set-TextNumberLike $ws "D39" "7.387"
Set-TextValue $ws "E39" "  +2.04%  "

Set-TextNumberLike $ws "D40" "0.5314"
Set-TextValue $ws "E40" "  +3.47%  "

Set-TextNumberLike $ws "D41" "0.1724"
Set-TextValue $ws "E41" "  +1.46%  "

Set-TextValue $ws "E42" "  -0.40%  "

Set-TextNumberLike $ws "D43" "2.299"
Set-TextValue $ws "E43" "  +18.87%  "

Set-TextNumberLike $ws "D44" "8.749"
Set-TextValue $ws "E44" "  +1.96%  "

Set-TextNumberLike $ws "D45" "0.5094"
Set-TextValue $ws "E45" "  +6.77%  "

Set-TextNumberLike $ws "D46" "10.68"
Set-TextValue $ws "E46" "  +0.45%  "

# Rows 47 & 48 were swapped (Quant <-> NEARProtocol) in the source feed.
Set-TextValue $ws "B47" "NEARProtocol"
Set-TextValue $ws "C47" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextNumberLike $ws "D47" "1.700"
Set-TextValue $ws "E47" "  +2.31%  "

Set-TextValue $ws "B48" "Quant"
Set-TextValue $ws "C48" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextNumberLike $ws "D48" "105.60"
Set-TextValue $ws "E48" "  -0.44%  "

Set-TextValue $ws "E49" "  +0.07%  "

Set-TextNumberLike $ws "D50" "0.06381"
Set-TextValue $ws "E50" "  +0.77%  "

Set-TextNumberLike $ws "D51" "0.9345"
Set-TextValue $ws "E51" "  +2.80%  "
